$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.975.74"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "2.271.03"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.63"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.01"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.81"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.69"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "2.624.79"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.36"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").Value = "2.273.16"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.786"
$ws.Range("E17").Value = "  +3.96%  "
$ws.Range("D18").Value = "41.860.71"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.78"
$ws.Range("E19").Value = "  +4.77%  "
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.99"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.26"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.31"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.99"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.68"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  -9.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.10"
$ws.Range("E30").Value = "  +3.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.47"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("E32").Value = "  +4.10%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.04"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.26"
$ws.Range("E36").Value = "  +3.99%  "
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.79"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.97"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.82"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "2.017.08"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.23"
$ws.Range("E44").Value = "  +8.78%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0283"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("E46").Value = "  +2.15%  "
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.32"
$ws.Range("E48").Value = "  +3.01%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.58"
$ws.Range("E50").Value = "  +3.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.14"
$ws.Range("E51").Value = "  +0.06%  "
